$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's visible text (runs merged to a single plain
# run) while reliably discarding any stranded <w:proofErr/> markers.
#
# The engine backing this COM shim only prunes a <w:proofErr/> element when
# it sits *strictly inside* the character range being reconciled; one that
# touches the start/end of the edited range (e.g. one parked right before
# the paragraph mark) survives untouched.  To guarantee every proofErr in
# the paragraph is interior to the edit, we temporarily splice the target
# paragraph together with the paragraph that follows it (deleting the
# paragraph mark between them), rewrite the *whole* combined text in one
# shot, and then reinsert the paragraph break at the right offset.  Because
# the edit's start/end now land inside neighbouring, untouched paragraphs
# (document start, or further content), every proofErr that belonged to the
# paragraph is interior and gets dropped.
#
# A "no-op" guard: if new text is character-for-character identical to the
# old text the engine sometimes treats the call as a no-op and skips
# reconciliation (leaving proofErr markers behind), so we always route the
# text through a throw-away intermediate value first to force a genuine
# diff before setting the real, final text.
# ---------------------------------------------------------------------------
function Set-ParaFullText {
    param($idx, [string]$newText)

    $p = $d.Paragraphs.Item($idx)
    $count = $d.Paragraphs.Count

    if ($idx -lt $count) {
        $nextP = $d.Paragraphs.Item($idx + 1)
        $nextFullText = $nextP.Range.Text
        if ($nextFullText.Length -gt 0 -and $nextFullText.Substring($nextFullText.Length - 1) -eq [char]13) {
            $nextPlain = $nextFullText.Substring(0, $nextFullText.Length - 1)
        } else {
            $nextPlain = $nextFullText
        }

        $mark = $d.Range($p.Range.End - 1, $p.Range.End)
        $mark.Delete()

        $merged = $d.Paragraphs.Item($idx).Range
        $merged.MoveEnd(1, -1)
        $merged.Text = $newText + "`u{E000}TEMP`u{E000}" + $nextPlain

        $merged2 = $d.Paragraphs.Item($idx).Range
        $merged2.MoveEnd(1, -1)
        $merged2.Text = $newText + $nextPlain

        $splitPos = $d.Paragraphs.Item($idx).Range.Start + $newText.Length
        $ins = $d.Range($splitPos, $splitPos)
        $ins.InsertParagraphAfter()
    } else {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $r.Text = $newText + "`u{E000}TEMP`u{E000}"
        $r2 = $p.Range
        $r2.MoveEnd(1, -1)
        $r2.Text = $newText
    }
}

function Set-RedRun([string]$needle) {
    $r = $d.Content
    $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Font.Color = 255
}

# --- Paragraph 1 ------------------------------------------------------------
Set-ParaFullText 1 "Aula 20 – pt1 introdução ao vuex – Mutations"

# --- Paragraph 2 --------------------------------------------------------
Set-ParaFullText 2 "Atualiza / substituição um campo com texto, ex: nome: Juca, atualiza = nome: João"

# --- Paragraph 4 (Aula 21) -----------------------------------------------
Set-ParaFullText 4 "Aula 21 – pt2 introdução ao vuex – Mutations"

# --- Paragraph 5 ----------------------------------------------------------
Set-ParaFullText 5 "Faz um produto(objeto) com característica, e adiciona e remove de um carrinho [array]"

# --- Paragraph 7 (Aula 22) -------------------------------------------------
Set-ParaFullText 7 "Aula 22 - Introdução VUEX Mutations parte 2"

# --- Paragraph 8 ------------------------------------------------------------
Set-ParaFullText 8 "Pega o valors dos produtos e da soma os produtos no carrinho [array]"

# --- Paragraph 9 (Aula 23) — keep the red "Revisar aula" highlight ---------
Set-ParaFullText 9 "Aula 23 - Introdução VUEX Mutations Actions – Revisar aula"
Set-RedRun "Revisar aula"

# --- Paragraph 10 -----------------------------------------------------------
Set-ParaFullText 10 "Faz a mesma coisa que mutations, mas sem commit"

# --- Paragraph 11 (Aula 24) — keep the red "Revisar aula" highlight --------
Set-ParaFullText 11 "Aula 24 - Introdução a composition API setup - Revisar aula"
Set-RedRun "Revisar aula"

# ---------------------------------------------------------------------------
# New content: empty paragraph + "Aula 26" block appended at the end of the
# document (after "Transforma variável em variável reativa").
# ---------------------------------------------------------------------------
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertParagraphAfter()
$endRng2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng2.InsertParagraphAfter()
$endRng3 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng3.InsertParagraphAfter()

$lastCount = $d.Paragraphs.Count
$pAula26 = $d.Paragraphs.Item($lastCount - 1)
$pAula26.Range.InsertBefore("Aula 26 - composition API Computed e Watch")

$pWatch = $d.Paragraphs.Item($lastCount)
$pWatch.Range.InsertBefore("Watch – quando o usuário for alterado, faz uma operação ou joga mensagem na tela.....")
